# Applies the "conclusions 3rd personed, and figures removed" edits to the
# ResponseDocument-style reviewer spreadsheet (Sheet1, columns A=Reviewer
# comments, B=Response, C=Changes to thesis).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) C4 - rich text run: "...restructured to 3^rd person." gains a new
#    trailing sentence. Append the new text, then re-stamp the run
#    formatting (Cambria 12, with the "rd" kept as a superscript) so the
#    saved shared string keeps its original look.
# ---------------------------------------------------------------------
$c4 = $ws.Range("C4")
$c4full = $c4.Characters().Text
$c4len = $c4full.Length
$c4ins = $c4.Characters($c4len + 1, 0)
$c4ins.Text = " This is performed throughout the thesis except for the Aims in the conclusion chapter"

$c4new = $c4.Characters().Text
$c4newLen = $c4new.Length

$c4run1 = $c4.Characters(1, 61)
$c4run1.Font.Name = "Cambria"
$c4run1.Font.Size = 12

$c4run2 = $c4.Characters(62, 2)
$c4run2.Font.Name = "Cambria"
$c4run2.Font.Size = 12
$c4run2.Font.Superscript = $true

$c4run3 = $c4.Characters(64, $c4newLen - 63)
$c4run3.Font.Name = "Cambria"
$c4run3.Font.Size = 12

# ---------------------------------------------------------------------
# 2) C6 - plain text: trailing "." replaced with
#    " from the abstract header."
# ---------------------------------------------------------------------
$c6 = $ws.Range("C6")
$c6old = $c6.Value2
$c6.Value = $c6old.Substring(0, $c6old.Length - 1) + " from the abstract header."

# ---------------------------------------------------------------------
# 3) B40 / C40 - reviewer response text rewritten/expanded
# ---------------------------------------------------------------------
$ws.Range("B40").Value = "Fair point, I like the idea of having bvoc emission modelling introduced in one spot. I have left the subsection in Chapter 2 that refers to how MEGAN is implemented within GEOS-Chem, and also several sentences in chapter 3 that are relevant to the surrounding text."
$ws.Range("C40").Value = "Added Section 1.3.2: Biogenic emissions modelling, drawing sentences from several sections in chapters 1 and 2, and references to the section are added in chapters 1 2 and 3. "

# ---------------------------------------------------------------------
# 4) Row 43 - new reviewer exchange about figures/tables in conclusions
#    chapter: B43 (response) and C43 (changes to thesis) newly populated.
# ---------------------------------------------------------------------
$ws.Range("B43").Value = "OK it seems fair to have no figures or new material in the conclusions, however I would like to keep the table as it is much clearer than a long list of numbers that I would have otherwise."
$ws.Range("C43").Value = "Figure 5.1 moved into Chapter 3 results as Figure 3.21, with some text added there, and some text replaced in the Conclusions chapter. Figure 5.2 removed, sentence now references original figure."

# ---------------------------------------------------------------------
# 5) Cosmetic follow-on from editing the sheet: selection moves to B41.
# ---------------------------------------------------------------------
$ws.Range("B41").Select() | Out-Null
